# Updates cryptos list data (Coin, Link, Price, Volume(1h)) for rows 2-51
# Commit: Updated cryptos list on Mon Dec  2 09:38:33 UTC 2024 with GitHub Actions
#
# All target columns (B/C/D/E) hold text in the original workbook (t="inlineStr").
# Force the NumberFormat to text ("@") before assigning so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00", "2.34") into numbers,
# which would silently drop formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '95.376.81'
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -1.62%  '

# Row 3
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.619.07'
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  -2.16%  '

# Row 4
$ws.Range("B4").NumberFormat = '@'
$ws.Range("B4").Value = 'XRP'
$ws.Range("C4").NumberFormat = '@'
$ws.Range("C4").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '2.34'
$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  +21.79%  '

# Row 5
$ws.Range("B5").NumberFormat = '@'
$ws.Range("B5").Value = 'TetherUSD'
$ws.Range("C5").NumberFormat = '@'
$ws.Range("C5").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '1.00'
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  +0.16%  '

# Row 6
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '226.94'
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  -4.48%  '

# Row 7
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '637.59'
$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  -2.87%  '

# Row 8
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.414'

# Row 9
$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  +2.58%  '

# Row 10
$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  +0.08%  '

# Row 11
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '3.618.62'
$ws.Range("E11").NumberFormat = '@'
$ws.Range("E11").Value = '  -2.11%  '

# Row 12
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '46.71'
$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  +5.56%  '

# Row 13
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  -0.83%  '

# Row 14
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  -2.27%  '

# Row 15
$ws.Range("E15").NumberFormat = '@'
$ws.Range("E15").Value = '  -4.15%  '

# Row 16
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '4.295.15'
$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  -2.14%  '

# Row 17
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '95.027.85'
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  -1.83%  '

# Row 18
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '8.78'
$ws.Range("E18").NumberFormat = '@'
$ws.Range("E18").Value = '  -1.78%  '

# Row 19
$ws.Range("B19").NumberFormat = '@'
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").NumberFormat = '@'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '19.91'
$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  +6.74%  '

# Row 20
$ws.Range("B20").NumberFormat = '@'
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").NumberFormat = '@'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '3.614.17'
$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  -2.74%  '

# Row 21
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '12.88'
$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  -1.08%  '

# Row 22
$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  +1.55%  '

# Row 23
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '511.96'
$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -2.13%  '

# Row 24
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '3.26'
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  -4.95%  '

# Row 25
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '0.245'
$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  +29.13%  '

# Row 26
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '119.90'
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  +17.78%  '

# Row 27
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  -3.30%  '

# Row 28
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '6.75'
$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  -2.56%  '

# Row 29
$ws.Range("B29").NumberFormat = '@'
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").NumberFormat = '@'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '3.814.40'
$ws.Range("E29").NumberFormat = '@'
$ws.Range("E29").Value = '  -2.14%  '

# Row 30
$ws.Range("B30").NumberFormat = '@'
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").NumberFormat = '@'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '12.69'
$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  -5.26%  '

# Row 31
$ws.Range("B31").NumberFormat = '@'
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").NumberFormat = '@'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '12.63'
$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  +2.50%  '

# Row 32
$ws.Range("B32").NumberFormat = '@'
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").NumberFormat = '@'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '2.92'
$ws.Range("E32").NumberFormat = '@'
$ws.Range("E32").Value = '  -2.98%  '

# Row 33
$ws.Range("B33").NumberFormat = '@'
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").NumberFormat = '@'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("B34").NumberFormat = '@'
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").NumberFormat = '@'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  +0.06%  '

# Row 35
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.78'
$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  -4.35%  '

# Row 36
$ws.Range("B36").NumberFormat = '@'
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").NumberFormat = '@'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.178'
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  -5.46%  '

# Row 37
$ws.Range("B37").NumberFormat = '@'
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").NumberFormat = '@'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '31.82'
$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  -1.31%  '

# Row 38
$ws.Range("B38").NumberFormat = '@'
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").NumberFormat = '@'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.587'
$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  -2.07%  '

# Row 39
$ws.Range("B39").NumberFormat = '@'
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").NumberFormat = '@'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("B40").NumberFormat = '@'
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").NumberFormat = '@'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '597.18'
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  -7.71%  '

# Row 41
$ws.Range("B41").NumberFormat = '@'
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").NumberFormat = '@'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '8.34'
$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  -5.66%  '

# Row 42
$ws.Range("B42").NumberFormat = '@'
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").NumberFormat = '@'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '6.82'
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -0.64%  '

# Row 43
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '40.38'
$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  -0.17%  '

# Row 44
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.159'
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  -0.89%  '

# Row 45
$ws.Range("B45").NumberFormat = '@'
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").NumberFormat = '@'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.481'
$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  +6.57%  '

# Row 46
$ws.Range("B46").NumberFormat = '@'
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").NumberFormat = '@'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0476'
$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  +3.52%  '

# Row 47
$ws.Range("B47").NumberFormat = '@'
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").NumberFormat = '@'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  -6.43%  '

# Row 48
$ws.Range("B48").NumberFormat = '@'
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").NumberFormat = '@'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.923'
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  -3.65%  '

# Row 49
$ws.Range("B49").NumberFormat = '@'
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").NumberFormat = '@'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '23.46'
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  -0.70%  '

# Row 50
$ws.Range("B50").NumberFormat = '@'
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").NumberFormat = '@'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '8.56'
$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  +0.08%  '

# Row 51
$ws.Range("B51").NumberFormat = '@'
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").NumberFormat = '@'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.21'
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -3.60%  '
